$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data right
$ws.Columns.Item(1).Insert()

# Set header and value for the new ID column
$ws.Range("A1").Value = "ID"
$ws.Range("A2").Value = 1

# Update selection to match target state
$ws.Range("B6").Select()
